$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the empty placeholder cells in column B (rows 2-10) so they are
# entirely removed from the sheet (not just emptied of visible content).
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Clear()
}

# Shift the fixed "mesa" position numbers in column A: rows 6-9 move up
# one (5->6, 6->7, 7->8, 8->9) and what used to be row 10's value (9)
# now becomes 5, placed on row 10.
$ws.Range("A6").Value = 6
$ws.Range("A7").Value = 7
$ws.Range("A8").Value = 8
$ws.Range("A9").Value = 9
$ws.Range("A10").Value = 5

# D10's capacity is now stored as a text value "3" instead of the numeric 4.
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3"
